# Aula 02 BD - update deck:
#  1) Fixed "updated automatically"-style date text "30/07/2021" -> "03/08/2021"
#     on the slide master, every slide layout, and the handout master.
#  2) Slide 12 body text "em 1984." -> "em 1988."

$p = $ppt.ActivePresentation

$oldDate = "30/07/2021"
$newDate = "03/08/2021"

function Update-DateShape($shp) {
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1a) Slide master date placeholder.
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

# 1b) Every slide layout's date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}

# 2) Slide 12 body copy: "em 1984." -> "em 1988."
$slide = $p.Slides.Item(12)
$shape = $slide.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange
$fullText = $tr.Text
$needle = "em 1984."
$idx = $fullText.IndexOf($needle)
if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $needle.Length)
    $sub.Text = "em 1988."
}
